{"js": "// Insert a new paragraph \"- In real life we would minify the code\"\n// immediately after the paragraph \"- single JS and CSS files\".\nconst body = context.document.body;\nconst results = body.search(\"- single JS and CSS files\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target paragraph text not found\");\n}\n\nconst targetRange = results.items[0];\nconst targetParagraph = targetRange.paragraphs.getFirst();\nconst newParagraph = targetParagraph.insertParagraph(\n  \"- In real life we would minify the code\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$range = $d.Content\n$range.Find.ClearFormatting()\n$found = $range.Find.Execute(\"- single JS and CSS files\")\n\nif ($found) {\n    $range.Collapse(0)  # wdCollapseEnd\n    $range.InsertParagraphAfter()\n    $range.Collapse(0)\n    $range.MoveStart(1, 1) | Out-Null\n    $range.Text = \"- In real life we would minify the code\"\n}\n"}
